# Mise à jour de l'application
# Adds the 6 new "J-1" training rows (435:440, date 2025-09-05 / serial 45905)
# that were recorded to Feuil1 under the existing "Entrainement / Global" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed rows 435:440 from the last existing row (434) so the new rows inherit
# its formatting (date style on B, centered style on D, etc.), then overwrite
# every cell with the real data below.
$ws.Range("A434:V434").Copy($ws.Range("A435:V440"))

$newRows = @(
    @{ Row = 435; A = "Entrainement"; B = 45905; C = "Global"; D = "J-1"; E = "Karahali Souaré";  F = "right forward";    G = "01:02:37"; H = 4.06; I = 0.12; J = 3.94; K = 0.11; L = 0.01; M = 0;    N = 0; O = 0; P = 3.51; Q = 22.79; R = 4.39; S = 19; T = 2; U = 21; V = 7 },
    @{ Row = 436; A = "Entrainement"; B = 45905; C = "Global"; D = "J-1"; E = "Jeremie Laurent";   F = "left forward";     G = "01:05:23"; H = 5.22; I = 0.33; J = 4.88; K = 0.26; L = 0.07; M = 0.01; N = 0; O = 1; P = 4.74; Q = 25.6;  R = 5.19; S = 28; T = 6; U = 26; V = 5 },
    @{ Row = 437; A = "Entrainement"; B = 45905; C = "Global"; D = "J-1"; E = "Ilyes Boughanmi";   F = "center forward";   G = "01:04:58"; H = 3.75; I = 0.03; J = 3.71; K = 0.03; L = 0;    M = 0;    N = 0; O = 0; P = 3.4;  Q = 18.66; R = 4.26; S = 18; T = 1; U = 5;  V = 3 },
    @{ Row = 438; A = "Entrainement"; B = 45905; C = "Global"; D = "J-1"; E = "Emmanuel Valey";    F = "left forward";     G = "01:05:25"; H = 4.87; I = 0.26; J = 4.6;  K = 0.22; L = 0.04; M = 0;    N = 0; O = 1; P = 4.27; Q = 25;    R = 5.12; S = 36; T = 4; U = 20; V = 3 },
    @{ Row = 439; A = "Entrainement"; B = 45905; C = "Global"; D = "J-1"; E = "Omar Benyounes";    F = "center midfield";  G = "01:04:26"; H = 4.61; I = 0.31; J = 4.29; K = 0.23; L = 0.08; M = 0.01; N = 0; O = 2; P = 4.16; Q = 25.31; R = 4.92; S = 27; T = 6; U = 23; V = 6 },
    @{ Row = 440; A = "Entrainement"; B = 45905; C = "Global"; D = "J-1"; E = "Mattheo Haon";      F = "right back";       G = "01:04:51"; H = 4.77; I = 0.24; J = 4.52; K = 0.2;  L = 0.05; M = 0;    N = 0; O = 0; P = 4.36; Q = 24.6;  R = 4.21; S = 20; T = 3; U = 10; V = 2 }
)

$textCols = @("A", "C", "D", "E", "F", "G")
$numCols  = @("H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V")

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    foreach ($col in $textCols) {
        $ws.Range("$col$r").Value = [string]$rowData[$col]
    }

    $ws.Range("B$r").Value = $rowData.B

    foreach ($col in $numCols) {
        $ws.Range("$col$r").Value = [double]$rowData[$col]
    }
}

# Reflect the post-entry selection/scroll state (Excel moved the cursor down
# to the next blank row after the newly pasted block).
$ws.Range("C445").Select()
